# Update the "想去人数" (want-to-go count) figures in column F for the
# two sheets that list exhibitions: "展览" and "全部类型". The rows are
# aligned by event name, but the row numbers differ slightly between the
# two sheets because "全部类型" also interleaves rows from the "演出" sheet.

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value map for the "展览" sheet.
$exhibitionRows = @{
    5  = 15739
    9  = 15457
    10 = 58
    11 = 9056
    12 = 386
    15 = 98
    16 = 201
    18 = 202
    20 = 59
    21 = 557
    22 = 26
    23 = 11
    26 = 4
    27 = 17
    29 = 92
    30 = 33
    34 = 42
    37 = 458
    39 = 5567
}

# Row -> new value map for the "全部类型" sheet (same events, some rows shifted).
$allTypesRows = @{
    5  = 15739
    9  = 15457
    10 = 58
    11 = 9056
    12 = 386
    15 = 98
    16 = 201
    18 = 202
    20 = 59
    21 = 557
    22 = 26
    23 = 11
    26 = 4
    27 = 17
    29 = 92
    30 = 33
    36 = 42
    39 = 458
    41 = 5567
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionRows.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionRows[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesRows.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesRows[$row]
}
